$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (Actual Outcome) for rows 2-6: "-" -> "Same as expected outcome."
$ws.Range("F2:F6").Value = "Same as expected outcome."

# Column G (Fail/Pass) for rows 2-6: "Fail" -> "Pass"
$ws.Range("G2:G6").Value = "Pass"

# Update the active selection to match the saved view state
$ws.Range("F2:F6").Select() | Out-Null
